# Adds the new Faker-generated "customer" fixture rows used by the
# "customer creation end to end test": 23 additional rows on the
# "registered" sheet (rows 12-34) and 1 additional row on the
# "jhi_user" sheet (row 10).

$wb = $excel.ActiveWorkbook
$wsRegistered = $wb.Worksheets.Item("registered")
$wsJhiUser = $wb.Worksheets.Item("jhi_user")

# --- "registered" sheet: rows 12-34 ---
# row 12
$wsRegistered.Range("A12").Value = '898-74-5291'
$wsRegistered.Range("B12").Value = 'Hang'
$wsRegistered.Range("C12").Value = 'Bartell'
$wsRegistered.Range("D12").Value = '293 Rufus Extensions'
$wsRegistered.Range("E12").NumberFormat = "@"
$wsRegistered.Range("E12").Value = '01532'
$wsRegistered.Range("E12").Style = "Normal"
$wsRegistered.Range("F12").Value = 'Swaniawskiside'
$wsRegistered.Range("G12").Value = 'Louisiana'
$wsRegistered.Range("H12").Value = '118-675-2045'
$wsRegistered.Range("I12").Value = '511-434-3547'
$wsRegistered.Range("J12").Value = 'marianela.funk'
$wsRegistered.Range("K12").Value = 'nicolas.zemlak@yahoo.com'
$wsRegistered.Range("L12").Value = 'uG2#I/^'

# row 13
$wsRegistered.Range("A13").Value = '309-75-3963'
$wsRegistered.Range("B13").Value = 'Maryellen'
$wsRegistered.Range("C13").Value = 'Kemmer'
$wsRegistered.Range("D13").Value = '7197 Thiel Well'
$wsRegistered.Range("E13").Value = '22601-3800'
$wsRegistered.Range("F13").Value = 'New Ashantifurt'
$wsRegistered.Range("G13").Value = 'Oklahoma'
$wsRegistered.Range("H13").Value = '967-954-5421'
$wsRegistered.Range("I13").Value = '118-675-2045'
$wsRegistered.Range("J13").Value = 'sigrid.hamill'
$wsRegistered.Range("K13").Value = 'cordia.corkery@gmail.com'
$wsRegistered.Range("L13").Value = 'eI6&7TB'

# row 14
$wsRegistered.Range("A14").Value = '836-33-2977'
$wsRegistered.Range("B14").Value = 'Allen'
$wsRegistered.Range("C14").Value = 'Steuber'
$wsRegistered.Range("D14").Value = '710 Rosanna Drive'
$wsRegistered.Range("E14").NumberFormat = "@"
$wsRegistered.Range("E14").Value = '14384'
$wsRegistered.Range("E14").Style = "Normal"
$wsRegistered.Range("F14").Value = 'North Nedra'
$wsRegistered.Range("G14").Value = 'Illinois'
$wsRegistered.Range("H14").Value = '511-434-3547'
$wsRegistered.Range("I14").Value = '319-959-3517'
$wsRegistered.Range("J14").Value = 'daren.cormier'
$wsRegistered.Range("K14").Value = 'eloy.stokes@yahoo.com'
$wsRegistered.Range("L14").Value = 'hW9–KAZ'

# row 15
$wsRegistered.Range("A15").Value = '754-05-5895'
$wsRegistered.Range("B15").Value = 'Randi'
$wsRegistered.Range("C15").Value = 'Lowe'
$wsRegistered.Range("D15").Value = '53389 Barton Orchard'
$wsRegistered.Range("E15").Value = '54470-5132'
$wsRegistered.Range("F15").Value = 'Strosinland'
$wsRegistered.Range("G15").Value = 'New Mexico'
$wsRegistered.Range("H15").Value = '189-843-0512'
$wsRegistered.Range("I15").Value = '715-457-6770'
$wsRegistered.Range("J15").Value = 'haley.treutel'
$wsRegistered.Range("K15").Value = 'lupe.funk@yahoo.com'
$wsRegistered.Range("L15").Value = 'dJ3(KUV'

# row 16
$wsRegistered.Range("A16").Value = '802-92-2872'
$wsRegistered.Range("B16").Value = 'Lelah'
$wsRegistered.Range("C16").Value = 'Waelchi'
$wsRegistered.Range("D16").Value = '0984 Haley Keys'
$wsRegistered.Range("E16").NumberFormat = "@"
$wsRegistered.Range("E16").Value = '38071'
$wsRegistered.Range("E16").Style = "Normal"
$wsRegistered.Range("F16").Value = 'West Jewell'
$wsRegistered.Range("G16").Value = 'Kansas'
$wsRegistered.Range("H16").Value = '189-843-0512'
$wsRegistered.Range("I16").Value = '967-954-5421'
$wsRegistered.Range("J16").Value = 'kendall.labadie'
$wsRegistered.Range("K16").Value = 'julio.connelly@gmail.com'
$wsRegistered.Range("L16").Value = 'gI0{n?*'

# row 17
$wsRegistered.Range("A17").Value = '357-68-1614'
$wsRegistered.Range("B17").Value = 'Asa'
$wsRegistered.Range("C17").Value = 'Schmidt'
$wsRegistered.Range("D17").Value = '1810 Anisa Street'
$wsRegistered.Range("E17").Value = '02909-4324'
$wsRegistered.Range("F17").Value = 'Kuhicview'
$wsRegistered.Range("G17").Value = 'Delaware'
$wsRegistered.Range("H17").Value = '588-697-5116'
$wsRegistered.Range("I17").Value = '588-697-5116'
$wsRegistered.Range("J17").Value = 'herbert.harris'
$wsRegistered.Range("K17").Value = 'rich.erdman@hotmail.com'
$wsRegistered.Range("L17").Value = 'tM3=<>('

# row 18
$wsRegistered.Range("A18").Value = '325-06-9910'
$wsRegistered.Range("B18").Value = 'Amber'
$wsRegistered.Range("C18").Value = 'Pacocha'
$wsRegistered.Range("D18").Value = '696 Nickolas Place'
$wsRegistered.Range("E18").NumberFormat = "@"
$wsRegistered.Range("E18").Value = '99297'
$wsRegistered.Range("E18").Style = "Normal"
$wsRegistered.Range("F18").Value = 'Devonburgh'
$wsRegistered.Range("G18").Value = 'Missouri'
$wsRegistered.Range("H18").Value = '967-954-5421'
$wsRegistered.Range("I18").Value = '588-697-5116'
$wsRegistered.Range("J18").Value = 'wesley.wiegand'
$wsRegistered.Range("K18").Value = 'randy.hoppe@hotmail.com'
$wsRegistered.Range("L18").Value = 'aR5[8LL'

# row 19
$wsRegistered.Range("A19").Value = '609-80-4637'
$wsRegistered.Range("B19").Value = 'Corene'
$wsRegistered.Range("C19").Value = 'Harber'
$wsRegistered.Range("D19").Value = '44609 Kuvalis Crossing'
$wsRegistered.Range("E19").Value = '03351-0751'
$wsRegistered.Range("F19").Value = 'Valorieland'
$wsRegistered.Range("G19").Value = 'Washington'
$wsRegistered.Range("H19").Value = '715-457-6770'
$wsRegistered.Range("I19").Value = '967-954-5421'
$wsRegistered.Range("J19").Value = 'michal.kohler'
$wsRegistered.Range("K19").Value = 'moses.satterfield@hotmail.com'
$wsRegistered.Range("L19").Value = 'xN0#S5v'

# row 20
$wsRegistered.Range("A20").Value = '851-01-8508'
$wsRegistered.Range("B20").Value = 'Jason'
$wsRegistered.Range("C20").Value = 'Heller'
$wsRegistered.Range("D20").Value = '02048 Prosacco Highway'
$wsRegistered.Range("E20").Value = '03137-3518'
$wsRegistered.Range("F20").Value = 'Port Garretmouth'
$wsRegistered.Range("G20").Value = 'Idaho'
$wsRegistered.Range("H20").Value = '107-834-0930'
$wsRegistered.Range("I20").Value = '107-834-0930'
$wsRegistered.Range("J20").Value = 'lillia.dickens'
$wsRegistered.Range("K20").Value = 'carey.monahan@hotmail.com'
$wsRegistered.Range("L20").Value = 'eB9(gcS'

# row 21
$wsRegistered.Range("A21").Value = '034-90-3565'
$wsRegistered.Range("B21").Value = 'Douglas'
$wsRegistered.Range("C21").Value = 'Becker'
$wsRegistered.Range("D21").Value = '38589 Landon Forest'
$wsRegistered.Range("E21").NumberFormat = "@"
$wsRegistered.Range("E21").Value = '69586'
$wsRegistered.Range("E21").Style = "Normal"
$wsRegistered.Range("F21").Value = 'Rickyburgh'
$wsRegistered.Range("G21").Value = 'Nevada'
$wsRegistered.Range("H21").Value = '319-959-3517'
$wsRegistered.Range("I21").Value = '107-834-0930'
$wsRegistered.Range("J21").Value = 'josef.gibson'
$wsRegistered.Range("K21").Value = 'cleopatra.steuber@gmail.com'
$wsRegistered.Range("L21").Value = 'hA0>5/f'

# row 22
$wsRegistered.Range("A22").Value = '336-17-2430'
$wsRegistered.Range("B22").Value = 'Suzy'
$wsRegistered.Range("C22").Value = 'Miller'
$wsRegistered.Range("D22").Value = '72337 Robert Trail'
$wsRegistered.Range("E22").Value = '55097-9016'
$wsRegistered.Range("F22").Value = 'Lake Chantaymouth'
$wsRegistered.Range("G22").Value = 'Illinois'
$wsRegistered.Range("H22").Value = '967-954-5421'
$wsRegistered.Range("I22").Value = '319-959-3517'
$wsRegistered.Range("J22").Value = 'vonda.schmidt'
$wsRegistered.Range("K22").Value = 'farrah.mante@gmail.com'
$wsRegistered.Range("L22").Value = 'iB7/Ura'

# row 23
$wsRegistered.Range("A23").Value = '054-41-1762'
$wsRegistered.Range("B23").Value = 'Marc'
$wsRegistered.Range("C23").Value = 'Considine'
$wsRegistered.Range("D23").Value = '402 Goldner Squares'
$wsRegistered.Range("E23").Value = '74219-0957'
$wsRegistered.Range("F23").Value = 'East Allenfurt'
$wsRegistered.Range("G23").Value = 'Oklahoma'
$wsRegistered.Range("H23").Value = '206-381-6032'
$wsRegistered.Range("I23").Value = '967-954-5421'
$wsRegistered.Range("J23").Value = 'debbra.effertz'
$wsRegistered.Range("K23").Value = 'cyndi.windler@hotmail.com'
$wsRegistered.Range("L23").Value = 'mE6}7bM'

# row 24
$wsRegistered.Range("A24").Value = '709-81-8304'
$wsRegistered.Range("B24").Value = 'Herta'
$wsRegistered.Range("C24").Value = 'Kuhic'
$wsRegistered.Range("D24").Value = '1976 Gino Mall'
$wsRegistered.Range("E24").Value = '33194-4593'
$wsRegistered.Range("F24").Value = 'West Randa'
$wsRegistered.Range("G24").Value = 'New Mexico'
$wsRegistered.Range("H24").Value = '189-843-0512'
$wsRegistered.Range("I24").Value = '319-959-3517'
$wsRegistered.Range("J24").Value = 'lashon.macgyver'
$wsRegistered.Range("K24").Value = 'ellis.bernier@hotmail.com'
$wsRegistered.Range("L24").Value = 'iB7)v^<'

# row 25
$wsRegistered.Range("A25").Value = '057-63-8250'
$wsRegistered.Range("B25").Value = 'Logan'
$wsRegistered.Range("C25").Value = 'Maggio'
$wsRegistered.Range("D25").Value = '2028 Parker Forges'
$wsRegistered.Range("E25").NumberFormat = "@"
$wsRegistered.Range("E25").Value = '08340'
$wsRegistered.Range("E25").Style = "Normal"
$wsRegistered.Range("F25").Value = 'Lake Jordon'
$wsRegistered.Range("G25").Value = 'Pennsylvania'
$wsRegistered.Range("H25").Value = '967-954-5421'
$wsRegistered.Range("I25").Value = '715-457-6770'
$wsRegistered.Range("J25").Value = 'reid.klocko'
$wsRegistered.Range("K25").Value = 'octavio.zboncak@gmail.com'
$wsRegistered.Range("L25").Value = 'wU4–;w!'

# row 26
$wsRegistered.Range("A26").Value = '859-83-4387'
$wsRegistered.Range("B26").Value = 'Raisa'
$wsRegistered.Range("C26").Value = 'Tillman'
$wsRegistered.Range("D26").Value = '076 Melvina Junctions'
$wsRegistered.Range("E26").Value = '46930-6203'
$wsRegistered.Range("F26").Value = 'North Quinn'
$wsRegistered.Range("G26").Value = 'New Mexico'
$wsRegistered.Range("H26").Value = '206-381-6032'
$wsRegistered.Range("I26").Value = '967-954-5421'
$wsRegistered.Range("J26").Value = 'andrea.farrell'
$wsRegistered.Range("K26").Value = 'reagan.dach@yahoo.com'
$wsRegistered.Range("L26").Value = 'sV8?2x8'

# row 27
$wsRegistered.Range("A27").Value = '825-87-7524'
$wsRegistered.Range("B27").Value = 'Lon'
$wsRegistered.Range("C27").Value = 'Morar'
$wsRegistered.Range("D27").Value = '861 Clifford Highway'
$wsRegistered.Range("E27").Value = '60521-6401'
$wsRegistered.Range("F27").Value = 'North Youlanda'
$wsRegistered.Range("G27").Value = 'West Virginia'
$wsRegistered.Range("H27").Value = '967-954-5421'
$wsRegistered.Range("I27").Value = '189-843-0512'
$wsRegistered.Range("J27").Value = 'lexie.heaney'
$wsRegistered.Range("K27").Value = 'reuben.vandervort@gmail.com'
$wsRegistered.Range("L27").Value = 'lG2>v(~'

# row 28
$wsRegistered.Range("A28").Value = '372-74-8288'
$wsRegistered.Range("B28").Value = 'Reyes'
$wsRegistered.Range("C28").Value = 'Rogahn'
$wsRegistered.Range("D28").Value = '8387 Rosenbaum Shoal'
$wsRegistered.Range("E28").NumberFormat = "@"
$wsRegistered.Range("E28").Value = '81643'
$wsRegistered.Range("E28").Style = "Normal"
$wsRegistered.Range("F28").Value = 'Mosciskimouth'
$wsRegistered.Range("G28").Value = 'Minnesota'
$wsRegistered.Range("H28").Value = '206-381-6032'
$wsRegistered.Range("I28").Value = '118-675-2045'
$wsRegistered.Range("J28").Value = 'marilyn.hegmann'
$wsRegistered.Range("K28").Value = 'raphael.purdy@yahoo.com'
$wsRegistered.Range("L28").Value = 'pG5}yrS'

# row 29
$wsRegistered.Range("A29").Value = '232-78-3042'
$wsRegistered.Range("B29").Value = 'Kerry'
$wsRegistered.Range("C29").Value = 'Pfannerstill'
$wsRegistered.Range("D29").Value = '231 Wunsch Brooks'
$wsRegistered.Range("E29").NumberFormat = "@"
$wsRegistered.Range("E29").Value = '61924'
$wsRegistered.Range("E29").Style = "Normal"
$wsRegistered.Range("F29").Value = 'East Chrisport'
$wsRegistered.Range("G29").Value = 'Virginia'
$wsRegistered.Range("H29").Value = '967-954-5421'
$wsRegistered.Range("I29").Value = '967-954-5421'
$wsRegistered.Range("J29").Value = 'omega.terry'
$wsRegistered.Range("K29").Value = 'matthew.boehm@yahoo.com'
$wsRegistered.Range("L29").Value = 'lJ4&c6?'

# row 30
$wsRegistered.Range("A30").Value = '304-63-9260'
$wsRegistered.Range("B30").Value = 'Nolan'
$wsRegistered.Range("C30").Value = 'Bartoletti'
$wsRegistered.Range("D30").Value = '3188 Beverley Ford'
$wsRegistered.Range("E30").Value = '20670-7308'
$wsRegistered.Range("F30").Value = 'New Aubrey'
$wsRegistered.Range("G30").Value = 'Ohio'
$wsRegistered.Range("H30").Value = '588-697-5116'
$wsRegistered.Range("I30").Value = '588-697-5116'
$wsRegistered.Range("J30").Value = 'jed.glover'
$wsRegistered.Range("K30").Value = 'tillie.aufderhar@gmail.com'
$wsRegistered.Range("L30").Value = 'qL8–Em)'

# row 31
$wsRegistered.Range("A31").Value = '500-26-8147'
$wsRegistered.Range("B31").Value = 'Marty'
$wsRegistered.Range("C31").Value = 'Daugherty'
$wsRegistered.Range("D31").Value = '5938 Treutel Grove'
$wsRegistered.Range("E31").NumberFormat = "@"
$wsRegistered.Range("E31").Value = '54061'
$wsRegistered.Range("E31").Style = "Normal"
$wsRegistered.Range("F31").Value = 'Lindview'
$wsRegistered.Range("G31").Value = 'Michigan'
$wsRegistered.Range("H31").Value = '118-675-2045'
$wsRegistered.Range("I31").Value = '967-954-5421'
$wsRegistered.Range("J31").Value = 'noelle.predovic'
$wsRegistered.Range("K31").Value = 'shaunna.schmitt@yahoo.com'
$wsRegistered.Range("L31").Value = 'zZ7<B{S'

# row 32
$wsRegistered.Range("A32").Value = '621-08-6725'
$wsRegistered.Range("B32").Value = 'Lilly'
$wsRegistered.Range("C32").Value = 'Brakus'
$wsRegistered.Range("D32").Value = '85879 Kovacek Roads'
$wsRegistered.Range("E32").NumberFormat = "@"
$wsRegistered.Range("E32").Value = '72758'
$wsRegistered.Range("E32").Style = "Normal"
$wsRegistered.Range("F32").Value = 'Anniechester'
$wsRegistered.Range("G32").Value = 'Indiana'
$wsRegistered.Range("H32").Value = '107-834-0930'
$wsRegistered.Range("I32").Value = '511-434-3547'
$wsRegistered.Range("J32").Value = 'christine.vonrueden'
$wsRegistered.Range("K32").Value = 'paris.nienow@gmail.com'
$wsRegistered.Range("L32").Value = 'bN6=ahw'

# row 33
$wsRegistered.Range("A33").Value = '490-65-9976'
$wsRegistered.Range("B33").Value = 'Ariel'
$wsRegistered.Range("C33").Value = 'Emmerich'
$wsRegistered.Range("D33").Value = '7320 Jenelle Plaza'
$wsRegistered.Range("E33").Value = '30643-4275'
$wsRegistered.Range("F33").Value = 'Turnerfurt'
$wsRegistered.Range("G33").Value = 'Iowa'
$wsRegistered.Range("H33").Value = '118-675-2045'
$wsRegistered.Range("I33").Value = '511-434-3547'
$wsRegistered.Range("J33").Value = 'erwin.lueilwitz'
$wsRegistered.Range("K33").Value = 'hugh.dare@hotmail.com'
$wsRegistered.Range("L33").Value = 'sU2=ZNa'

# row 34
$wsRegistered.Range("A34").Value = '201-11-9899'
$wsRegistered.Range("B34").Value = 'Majorie'
$wsRegistered.Range("C34").Value = 'DuBuque'
$wsRegistered.Range("D34").Value = '70675 Jacquelin Fields'
$wsRegistered.Range("E34").NumberFormat = "@"
$wsRegistered.Range("E34").Value = '06206'
$wsRegistered.Range("E34").Style = "Normal"
$wsRegistered.Range("F34").Value = 'West Kelsieborough'
$wsRegistered.Range("G34").Value = 'Kentucky'
$wsRegistered.Range("H34").Value = '715-457-6770'
$wsRegistered.Range("I34").Value = '319-959-3517'
$wsRegistered.Range("J34").Value = 'hang.crooks'
$wsRegistered.Range("K34").Value = 'lilla.fritsch@hotmail.com'
$wsRegistered.Range("L34").Value = 'nW4(nqC'

# --- "jhi_user" sheet: row 10 ---
# row 10
$wsJhiUser.Range("A10").Value = '584-05-9172'
$wsJhiUser.Range("B10").Value = 'Lewis'
$wsJhiUser.Range("C10").Value = 'Reichel'
$wsJhiUser.Range("D10").Value = '0495 Hilario Shores'
$wsJhiUser.Range("E10").Value = '15757-3265'
$wsJhiUser.Range("F10").Font.Bold = $false
$wsJhiUser.Range("G10").Value = 'Lake Enochland'
$wsJhiUser.Range("H10").Value = 'Arizona'
$wsJhiUser.Range("I10").Value = '967-954-5421'
$wsJhiUser.Range("J10").Font.Bold = $false
$wsJhiUser.Range("K10").Value = 'tommie.ritchie'
$wsJhiUser.Range("L10").Value = 'jesusa.tillman@yahoo.com'
$wsJhiUser.Range("M10").Value = 'wA1;U;i'
$wsJhiUser.Range("N10").Value = 'true'
$wsJhiUser.Range("O10").Value = 'ROLE_CUSTOMER'
$wsJhiUser.Range("P10").NumberFormat = "@"
$wsJhiUser.Range("P10").Value = '63465'
$wsJhiUser.Range("P10").Style = "Normal"

